{"js": "// Replace the date line and every two-digit-by-two-digit multiplication\n// answer in the practice table with the new values from the target\n// revision. Every \"old\" string below is unique within the document, so a\n// simple search-and-replace (first match) is unambiguous and safe.\nconst replacements = [\n  [\"2024-03-30 Saturday\", \"2024-03-31 Sunday\"],\n  [\"13\u00d726=338\", \"65\u00d723=1495\"],\n  [\"97\u00d714=1358\", \"88\u00d726=2288\"],\n  [\"66\u00d715=990\", \"69\u00d729=2001\"],\n  [\"74\u00d717=1258\", \"63\u00d795=5985\"],\n  [\"92\u00d750=4600\", \"81\u00d753=4293\"],\n  [\"30\u00d712=360\", \"19\u00d753=1007\"],\n  [\"43\u00d746=1978\", \"80\u00d775=6000\"],\n  [\"44\u00d777=3388\", \"36\u00d761=2196\"],\n  [\"77\u00d730=2310\", \"63\u00d767=4221\"],\n  [\"42\u00d720=840\", \"35\u00d795=3325\"],\n  [\"55\u00d733=1815\", \"51\u00d719=969\"],\n  [\"90\u00d745=4050\", \"58\u00d716=928\"],\n  [\"41\u00d741=1681\", \"55\u00d753=2915\"],\n  [\"89\u00d769=6141\", \"11\u00d789=979\"],\n  [\"95\u00d713=1235\", \"39\u00d777=3003\"],\n  [\"56\u00d756=3136\", \"73\u00d737=2701\"],\n  [\"31\u00d741=1271\", \"16\u00d756=896\"],\n  [\"69\u00d790=6210\", \"44\u00d797=4268\"],\n  [\"15\u00d736=540\", \"66\u00d772=4752\"],\n  [\"54\u00d759=3186\", \"91\u00d713=1183\"],\n  [\"63\u00d784=5292\", \"13\u00d739=507\"],\n  [\"40\u00d757=2280\", \"84\u00d775=6300\"],\n  [\"61\u00d740=2440\", \"40\u00d718=720\"],\n  [\"96\u00d744=4224\", \"65\u00d797=6305\"],\n  [\"71\u00d778=5538\", \"26\u00d723=598\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Could not find text: ${oldText}`);\n  }\n\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Replace the date line and every two-digit-by-two-digit multiplication\n# answer in the practice table with the new values from the target\n# revision. Every \"old\" string is unique within the document, so a plain\n# Find/Replace (first match, whole document range) is unambiguous.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"2024-03-30 Saturday\", \"2024-03-31 Sunday\"),\n  @(\"13\u00d726=338\", \"65\u00d723=1495\"),\n  @(\"97\u00d714=1358\", \"88\u00d726=2288\"),\n  @(\"66\u00d715=990\", \"69\u00d729=2001\"),\n  @(\"74\u00d717=1258\", \"63\u00d795=5985\"),\n  @(\"92\u00d750=4600\", \"81\u00d753=4293\"),\n  @(\"30\u00d712=360\", \"19\u00d753=1007\"),\n  @(\"43\u00d746=1978\", \"80\u00d775=6000\"),\n  @(\"44\u00d777=3388\", \"36\u00d761=2196\"),\n  @(\"77\u00d730=2310\", \"63\u00d767=4221\"),\n  @(\"42\u00d720=840\", \"35\u00d795=3325\"),\n  @(\"55\u00d733=1815\", \"51\u00d719=969\"),\n  @(\"90\u00d745=4050\", \"58\u00d716=928\"),\n  @(\"41\u00d741=1681\", \"55\u00d753=2915\"),\n  @(\"89\u00d769=6141\", \"11\u00d789=979\"),\n  @(\"95\u00d713=1235\", \"39\u00d777=3003\"),\n  @(\"56\u00d756=3136\", \"73\u00d737=2701\"),\n  @(\"31\u00d741=1271\", \"16\u00d756=896\"),\n  @(\"69\u00d790=6210\", \"44\u00d797=4268\"),\n  @(\"15\u00d736=540\", \"66\u00d772=4752\"),\n  @(\"54\u00d759=3186\", \"91\u00d713=1183\"),\n  @(\"63\u00d784=5292\", \"13\u00d739=507\"),\n  @(\"40\u00d757=2280\", \"84\u00d775=6300\"),\n  @(\"61\u00d740=2440\", \"40\u00d718=720\"),\n  @(\"96\u00d744=4224\", \"65\u00d797=6305\"),\n  @(\"71\u00d778=5538\", \"26\u00d723=598\")\n)\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
